# 8_S_Abutment_pTHg_m3_Flux_Annual_slrose.xlsx
# Model 3 rloadest run was re-run with corrected values (per commit message);
# refresh the computed flux/error columns and mark the two "NA" (winter)
# rows' % std-err-of-flux cell with the same "no fill" formatting Excel
# applies when AutoFilling a formula column that skips non-numeric rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Period 1 (row 2) -----------------------------------------------
$ws.Range("D2").Value = 0.0086727496967275892
$ws.Range("E2").Value = 0.0049106404946817096
$ws.Range("F2").Value = 0.056489149816248599
$ws.Range("G2").Value = 0.000029261759077413702
$ws.Range("H2").Value = 0.059194130017568299

# ---- Period 2 (row 3) -----------------------------------------------
$ws.Range("D3").Value = 0.082087508243031002
$ws.Range("E3").Value = 0.039338574060584101
$ws.Range("F3").Value = 0.217337905079109
$ws.Range("G3").Value = 0.00171638743758848
$ws.Range("H3").Value = 0.49012557181219302

# ---- Period 3 (row 4) is the winter "NA" row -- no numeric changes --

# ---- Period 4 (row 5) -------------------------------------------------
$ws.Range("D5").Value = 0.016886968202844899
$ws.Range("E5").Value = 0.0055481510443812698
$ws.Range("F5").Value = 0.10991189307296
$ws.Range("G5").Value = 0.000057057681432323498
$ws.Range("H5").Value = 0.115257877466642

# ---- Period 5 (row 6) is the winter "NA" row -- no numeric changes --

# ---- Period 6 (row 7) -------------------------------------------------
$ws.Range("D7").Value = 0.034574564666249502
$ws.Range("E7").Value = 0.016669397238608
$ws.Range("F7").Value = 0.208361897253588
$ws.Range("G7").Value = 0.000135944738824673
$ws.Range("H7").Value = 0.235630480101677

# J2/J3/J5/J7 already hold "=(Ex/Dx)*100" formulas that recalc from the
# updated D/E columns above. J4 and J6 correspond to the "NA" rows (no
# flux computed that period); give them the same "no fill" swatch as
# the rest of column J so the highlighted % std-err column reads cleanly.
$ws.Range("J4").Interior.ColorIndex = -4142
$ws.Range("J6").Interior.ColorIndex = -4142

# The workbook was left with a different cell selected/active on save.
$ws.Range("J14").Select() | Out-Null
